$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: add a new import line -------------------------------
# Insert a row INSIDE the existing D3:D4 merge (at row 4). Excel auto-
# extends the merge to D3:D5 and shifts the old row-4 content (blank D /
# "org.openl.rules.mapping.to" E) down to row 5, inheriting formatting.
$ws.Rows("4:4").Insert()

# Tidy up the formatting of the freshly inserted row 4 to match the rest
# of the merged label block / the E-column import rows.
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E5").Copy()
$ws.Range("E4").PasteSpecial(-4122)  # xlPasteFormats

# Row 4 keeps the original import text; row 5 (which now holds a shifted
# copy of that content) becomes the new import line required for the
# deep-mapping-with-array-in-the-middle scenario.
$ws.Range("E4").Value = "org.openl.rules.mapping.to"
$ws.Range("E5").Value = "org.openl.rules.mapping.to.containers"

# Add the second blank row (pushes everything from the old row 5 onward
# down by one more, so the old row 8 block lands on row 10).
$ws.Rows("6:6").Insert()
$ws.Rows("6:6").Clear()

# --- Data rows: add a parallel shallow array-index mapping -------------
# Rows 13-15 (old 11-13) gain a second destination column alongside the
# existing deep "array[N].b.aString" one: ArrayOfTypeBContainer's
# shallower "array[N].aString", exercising an index operator that sits in
# the middle of the destination path.
$ws.Range("E13").Value = "ArrayOfTypeBContainer"
$ws.Range("G13").Value = "array[1].aString"

$ws.Range("E14").Value = "ArrayOfTypeBContainer"
$ws.Range("G14").Value = "array[2].aString"

$ws.Range("E15").Value = "ArrayOfTypeBContainer"
$ws.Range("G15").Value = "array[3].aString"

# --- Selection moved to the newly-edited header block -------------------
$ws.Range("D2:E5").Select()
$excel.ActiveCell = $ws.Range("D2")
